$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40; rows 40-61 shift down to 41-62.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new weekly record.
$ws.Range("A40").Value = 5
$ws.Range("B40").Value = "Macroferia Regional de Talca"
$ws.Range("C40").Value = "Maule"
$ws.Range("D40").Value = [DateTime]"2022-02-04"
$ws.Range("E40").Value = 7
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100101
$ws.Range("H40").Value = "Berries"
$ws.Range("I40").Value = 100101001
$ws.Range("J40").Value = "Arándano (blue)"
$ws.Range("K40").Value = "Sin especificar"
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 200
$ws.Range("N40").Value = 3500
$ws.Range("O40").Value = 3500
$ws.Range("P40").Value = 3500
$ws.Range("Q40").Value = "$/bandeja 2 kilos"
$ws.Range("R40").Value = "Provincia de Linares"
$ws.Range("S40").Value = 1750
$ws.Range("T40").Value = 2
